# Update the "views" column (F) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "2"  = 61
    "3"  = 789
    "4"  = 41
    "6"  = 76
    "7"  = 275
    "8"  = 3959
    "10" = 4657
    "11" = 515
    "12" = 1174
    "13" = 76
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    "2"  = 61
    "3"  = 789
    "4"  = 41
    "6"  = 76
    "8"  = 275
    "9"  = 3959
    "11" = 4657
    "12" = 515
    "13" = 1174
    "14" = 76
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
